# 茶园面积 (tea garden area) data refresh:
#   - drop the 2000-2009 rows, shift 2010-2020 up to rows 2-12
#   - append a new 2021 row (no "本年采摘茶园面积" figure yet, same as 2020's
#     still-missing value)
#   - correct the 2019 "实有茶园面积" figure from 3104.8 to 3104.77

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the ten obsolete years (2000年..2009年, rows 2-11). This shifts the
# 2010年..2020年 rows (12-22) up to rows 2-12, carrying their formatting and
# values (including the still-blank 本年采摘茶园面积 cell on the old 2020年
# row) along with them.
$ws.Range("A2:C11").EntireRow.Delete()

# The 2019年 row (now row 11) gets a corrected 实有茶园面积 value.
$ws.Cells.Item(11, 2).Value = 3104.77

# Duplicate the now-last row (2020年, row 12) down into a fresh row 13 so the
# new 2021年 row inherits the same formatting (bold/centered/bordered year
# label style).
$ws.Range("A12:C12").Copy($ws.Range("A13:C13"))

# Overwrite the copied values with the new 2021年 figures; 本年采摘茶园面积
# isn't available yet, so that cell is left blank like 2020年's was.
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 3307.84106602
